# Apply the "pandoc title block" restructuring to the Day After Day article:
#   * drop the bookmark that used to wrap the old Heading1 title paragraph
#   * turn the title paragraph into a "Title"-styled paragraph, word-by-word runs
#   * turn "By Dorothy Day" into an "Authors"-styled "Dorothy Day" paragraph,
#     word-by-word runs (no more leading "By " / no more manual bold)

$d = $word.ActiveDocument

# --- 1. Best-effort removal of the (malformed) bookmark that brackets the
#        title paragraph. Some producers emit word-ish bookmark ids; try a
#        few ways of reaching it so this still works if the id is numeric.
$bookmarkName = "day-after-day---july-august-1935"
try {
    $bmCount = $d.Bookmarks.Count
    for ($i = $bmCount; $i -ge 1; $i--) {
        try {
            $bm = $d.Bookmarks.Item($i)
            if ($bm.Name -eq $bookmarkName) {
                $bm.Delete()
            }
        } catch {
        }
    }
} catch {
}
try {
    if ($d.Bookmarks.Exists($bookmarkName)) {
        $d.Bookmarks($bookmarkName).Delete()
    }
} catch {
}

# --- 2. Rebuild the title paragraph + the author paragraph in one shot.
#        Paragraph 1 is currently the Heading1 "Day After Day - July-August
#        1935" paragraph; paragraph 2 is the bold "By Dorothy Day" paragraph.
#        Replace both of them (text + paragraph marks) with freshly
#        constructed paragraphs, each run holding a single word/separator,
#        matching the pandoc-title-block export shape.

$titleWords   = @("Day", " ", "After", " ", "Day", " ", "-", " ", "July", "-", "August", " ", "1935")
$authorWords  = @("Dorothy", " ", "Day")

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-RunsXml($words) {
    $sb = ""
    foreach ($word in $words) {
        $escaped = $word.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $sb += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }
    return $sb
}

$titleXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + (New-RunsXml $titleWords) + '</w:p>'
$authorsXml = '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + (New-RunsXml $authorWords) + '</w:p>'

$firstPara  = $d.Paragraphs.Item(1)
$secondPara = $d.Paragraphs.Item(2)
$targetRange = $d.Range($firstPara.Range.Start, $secondPara.Range.End)
$targetRange.InsertXML($titleXml + $authorsXml)
